$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submenu")

$data = @(
    @("建物3", "Orange", 13),
    @("建物5", "Tomato", 17),
    @("建物6", "Yellow", 20),
    @("建物222", "Green", 8),
    @("建物42412", "IndianRed", 9),
    @("建物3434", "LightBlue", 10),
    @("建物31231", "Navy", 12),
    @("建物231", "Cyan", 7),
    @("建物434", "Salmon", 14)
)

$row = 5
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

$ws.Range("A8:XFD12").Select()

$excel.ActiveWindow.xWindow = 7650
$excel.ActiveWindow.yWindow = 5535
